# Applies the COVAC_TRACKER v1.1.2 / DHIS2.35.3 metadata refresh described
# in the commit "feat: update COVAC_tracker for DHIS2v 2.33, 2.34, 2.35".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Package info" sheet - version / build metadata bump
# ---------------------------------------------------------------------------
$pkg = $wb.Worksheets.Item("Package info")
$pkg.Range("B4").Value = "V1.1.2"
$pkg.Range("B5").Value = "DHIS2.35.3-3492688"
$pkg.Range("B6").Value = "20210408T081801"
$pkg.Range("B7").Value = "COVAC_TRACKER_V1.1.2_DHIS2.35.3-3492688_20210408T081801"
$pkg.Columns.Item(2).ColumnWidth = 57.7109375

# ---------------------------------------------------------------------------
# 2. "programTrackedEntityAttributes" sheet - replace raw UIDs with the
#    human readable attribute names in column B
# ---------------------------------------------------------------------------
$ptea = $wb.Worksheets.Item("programTrackedEntityAttributes")
$ptea.Range("B4").Value = "First Name"
$ptea.Range("B5").Value = "Surname"
$ptea.Range("B6").Value = "Sex"
$ptea.Range("B8").Value = "Date of birth"
$ptea.Range("B10").Value = "Home Address"

# ---------------------------------------------------------------------------
# 3. "dataElementGroups" sheet - the Data Element column got re-shuffled
#    (same 22 values, new row order)
# ---------------------------------------------------------------------------
$deg = $wb.Worksheets.Item("dataElementGroups")
$degValues = @(
    "COVAC - Underlying condition Other",
    "COVAC - Dose Number",
    "COVAC- Batch Number",
    "COVAC - Renal Disease",
    "COVAC - Malignancy",
    "COVAC - Vaccine Name",
    "COVAC - Pregnancy",
    "COVAC - Immunodeficiency",
    "COVAC Previously infected with COVID",
    "COVAC Suggested date for next dose",
    "COVAC - Multiple products used - Explain",
    "COVAC - AEFIs present",
    "COVAC - Cardiovascular Disease",
    "COVAC - Vaccine Manufacturer",
    "COVAC - Last Dose",
    "COVAC - Pregnancy gestation",
    "COVAC - Allergic reaction after first dose",
    "COVAC - Chronic Lung Disease",
    "COVAC - Diabetes",
    "COVAC - Neurological/Neuromuscular",
    "COVAC - Underlying condition",
    "COVAC - Total doses"
)
for ($i = 0; $i -lt $degValues.Length; $i++) {
    $row = 3 + $i
    $deg.Cells.Item($row, 2).Value = $degValues[$i]
}

# ---------------------------------------------------------------------------
# 4. "trackedEntityAttributes" sheet - attributes reordered and two new
#    rows inserted (Home Address, Surname) to bring the list in sync with
#    the tracked entity type attributes used by the program
# ---------------------------------------------------------------------------
$tea = $wb.Worksheets.Item("trackedEntityAttributes")

# Shift columns 4-12 downward isn't needed: rebuild rows 4-12 directly since
# the whole block changes (2 rows are new, others keep going but reordered).
$tea.Cells.Item(4, 1).Value = "Date of birth"
$tea.Cells.Item(4, 2).Value = "patinfo_ageonsetunit"
$tea.Cells.Item(4, 3).Value = ""
$tea.Cells.Item(4, 4).Value = ""
$tea.Cells.Item(4, 5).Value = "NI0QRzJvQ0k"

$tea.Cells.Item(5, 1).Value = "Date of birth is estimated"
$tea.Cells.Item(5, 2).Value = ""
$tea.Cells.Item(5, 3).Value = ""
$tea.Cells.Item(5, 4).Value = ""
$tea.Cells.Item(5, 5).Value = "Z1rLc1rVHK8"

$tea.Cells.Item(6, 1).Value = "First Name"
$tea.Cells.Item(6, 2).Value = "first_name"
$tea.Cells.Item(6, 3).Value = ""
$tea.Cells.Item(6, 4).Value = ""
$tea.Cells.Item(6, 5).Value = "sB1IHYu2xQT"

$tea.Cells.Item(7, 1).Value = "Home Address"
$tea.Cells.Item(7, 2).Value = "patinfo_resadmin0"
$tea.Cells.Item(7, 3).Value = ""
$tea.Cells.Item(7, 4).Value = ""
$tea.Cells.Item(7, 5).Value = "Xhdn49gUd52"

$tea.Cells.Item(8, 1).Value = "Mobile phone number"
$tea.Cells.Item(8, 2).Value = ""
$tea.Cells.Item(8, 3).Value = ""
$tea.Cells.Item(8, 4).Value = ""
$tea.Cells.Item(8, 5).Value = "fctSQp5nAYl"

$tea.Cells.Item(9, 1).Value = "National ID"
$tea.Cells.Item(9, 2).Value = ""
$tea.Cells.Item(9, 3).Value = ""
$tea.Cells.Item(9, 4).Value = ""
$tea.Cells.Item(9, 5).Value = "Ewi7FUfcHAD"

$tea.Cells.Item(10, 1).Value = "Sex"
$tea.Cells.Item(10, 2).Value = "patinfo_sex"
$tea.Cells.Item(10, 3).Value = ""
$tea.Cells.Item(10, 4).Value = ""
$tea.Cells.Item(10, 5).Value = "oindugucx72"

$tea.Cells.Item(11, 1).Value = "Surname"
$tea.Cells.Item(11, 2).Value = "surname"
$tea.Cells.Item(11, 3).Value = "The patient's surname (family name)"
$tea.Cells.Item(11, 4).Value = ""
$tea.Cells.Item(11, 5).Value = "ENRjVGxVL6l"

$tea.Cells.Item(12, 1).Value = "Unique System Identifier (EPI)"
$tea.Cells.Item(12, 2).Value = ""
$tea.Cells.Item(12, 3).Value = "System-generated unique ID following pattern: EPI prefix + value randomly generated (#####) - Customize the length depending on the target population of your implementation"
$tea.Cells.Item(12, 4).Value = ""
$tea.Cells.Item(12, 5).Value = "KSr2yTdu1AI"

$tea.Columns.Item(2).ColumnWidth = 22.7109375

# ---------------------------------------------------------------------------
# 5. "programs" sheet - bump the "Last updated" date
#    (force text format first so Excel doesn't reinterpret the literal as a
#    date serial number - the source stores it as a plain string)
# ---------------------------------------------------------------------------
$programs = $wb.Worksheets.Item("programs")
$programs.Range("C2").NumberFormat = "@"
$programs.Range("C2").Value = "2021-03-19"
